$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A167").Value = "FOS_ND200"
$ws.Range("B167").Value = "ABX_DISK"
$ws.Range("A168").Value = "FOS_NM"
$ws.Range("B168").Value = "ABX_MIC"

$ws.Range("D170").Select()
